$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a brand-new column at A; everything that was in A:X slides right
# to become B:Y (formulas/ranges/dataValidation/cols all shift automatically).
$ws.Columns("A:A").Insert()

# The new column A has no formatting yet - clone it from column B (which
# used to be column A) so header/data styling (borders, fill, alignment,
# bold) stays consistent with the rest of the sheet.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New index column content.
$ws.Range("A1").Value = "INDEX (DO NOT MODIFY)"
$ws.Range("A2").Value = 2

# Give the new index column its own width.
$ws.Columns("A:A").ColumnWidth = 22.17

# Uppercase every shifted header label except the last one (the "Status as
# of ..." column, which is left untouched).
$ws.Range("B1").Value = $ws.Range("B1").Text.ToUpper()
$ws.Range("C1").Value = $ws.Range("C1").Text.ToUpper()
$ws.Range("D1").Value = $ws.Range("D1").Text.ToUpper()
$ws.Range("E1").Value = $ws.Range("E1").Text.ToUpper()
$ws.Range("F1").Value = $ws.Range("F1").Text.ToUpper()
$ws.Range("G1").Value = $ws.Range("G1").Text.ToUpper()
$ws.Range("H1").Value = $ws.Range("H1").Text.ToUpper()
$ws.Range("I1").Value = $ws.Range("I1").Text.ToUpper()
$ws.Range("J1").Value = $ws.Range("J1").Text.ToUpper()
$ws.Range("K1").Value = $ws.Range("K1").Text.ToUpper()
$ws.Range("L1").Value = $ws.Range("L1").Text.ToUpper()
$ws.Range("M1").Value = $ws.Range("M1").Text.ToUpper()
$ws.Range("N1").Value = $ws.Range("N1").Text.ToUpper()
$ws.Range("O1").Value = $ws.Range("O1").Text.ToUpper()
$ws.Range("P1").Value = $ws.Range("P1").Text.ToUpper()
$ws.Range("Q1").Value = $ws.Range("Q1").Text.ToUpper()
$ws.Range("R1").Value = $ws.Range("R1").Text.ToUpper()
$ws.Range("S1").Value = $ws.Range("S1").Text.ToUpper()
$ws.Range("T1").Value = $ws.Range("T1").Text.ToUpper()
$ws.Range("U1").Value = $ws.Range("U1").Text.ToUpper()
$ws.Range("V1").Value = $ws.Range("V1").Text.ToUpper()
$ws.Range("W1").Value = $ws.Range("W1").Text.ToUpper()
$ws.Range("X1").Value = $ws.Range("X1").Text.ToUpper()
# Y1 ("Status as of July 11, 2025") keeps its original casing.
